# "Update roaming settings sample"
#
# The Snippets table gains one new data row: a "remove" member entry for the
# RoamingSettings class, inserted alphabetically between the existing "get"
# (row 288) and "saveAsync" (old row 289) entries. Every row from the old 289
# through 306 shifts down by one (to 290-307), and the table/autofilter/
# dimension all grow from F306 to F307.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 289; everything at/after 289 (old saveAsync.. row)
# shifts down to 290..307.
$ws.Cells.Item(289, 1).EntireRow.Insert()

# Grow the "Snippets" table/autofilter range so it covers the new row too.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F307"))

# Populate the new row: Office | RoamingSettings | remove | 1 |
#   outlook-roaming-settings-roaming-settings | remove
$ws.Cells.Item(289, 1).Value = "Office"

$ws.Cells.Item(289, 2).Value = "RoamingSettings"
$ws.Cells.Item(289, 2).Style = "Normal"

$ws.Cells.Item(289, 3).Value = "remove"
$ws.Cells.Item(289, 3).NumberFormat = "General"

$ws.Cells.Item(289, 4).Value = 1
$ws.Cells.Item(289, 4).NumberFormat = "General"
$ws.Cells.Item(289, 4).HorizontalAlignment = -4152 # xlRight

$ws.Cells.Item(289, 5).Value = "outlook-roaming-settings-roaming-settings"
$ws.Cells.Item(289, 5).Style = "Normal"

$ws.Cells.Item(289, 6).Value = "remove"
$ws.Cells.Item(289, 6).NumberFormat = "General"

# Leave the sheet scrolled/selected where the edit happened.
$ws.Activate()
$ws.Range("F289").Select()
